# Update column F (dSF) values for specific rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -2
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = -4
$ws.Range("F15").Value = -4
$ws.Range("F19").Value = -7
$ws.Range("F35").Value = -1
$ws.Range("F41").Value = 0
